# Added data driven test case and updated conftest to supress save password popup

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "login_data_sheet"

# Update the data driven test rows (typos introduced in the new test data)
$ws.Range("A3").Value = "adm@yourtore.com"
$ws.Range("A4").Value = "admin@youtore.com"

# Move active selection to H9
$ws.Range("H9").Select() | Out-Null
